$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1708
$ws.Range("J32").Value = 1708
$ws.Range("L32").Value = 1708
$ws.Range("N32").Value = -2360
$ws.Range("H64").Value = 4572.154
$ws.Range("J64").Value = 4356.25
$ws.Range("L64").Value = 4356.25
$ws.Range("N64").Value = -4852.25
$ws.Range("H67").Value = 4572.154
$ws.Range("J67").Value = 4356.25
$ws.Range("L67").Value = 4356.25
$ws.Range("N67").Value = -6072.25
$ws.Range("H100").Value = 4823.8945
$ws.Range("I100").Value = 3192
$ws.Range("K100").Value = 3192
$ws.Range("M100").Value = -2651
$ws.Range("H111").Value = 1052
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1052
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 3156
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -9290
$ws.Range("H127").Value = 3672.25
$ws.Range("I127").Value = 1466.6666
$ws.Range("J127").Value = 4995.6
$ws.Range("K127").Value = 4399.9998
$ws.Range("L127").Value = 14986.8
$ws.Range("M127").Value = 560.0002000000004
$ws.Range("N127").Value = -24906.8
$ws.Range("H129").Value = 2293.5625
$ws.Range("I129").Value = 983.5454999999999
$ws.Range("K129").Value = 2950.6365
$ws.Range("M129").Value = 2049.3635
$ws.Range("H137").Value = 19003580
$ws.Range("I137").Value = 33336056
$ws.Range("K137").Value = 100008168
$ws.Range("M137").Value = -100005618

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4959.8
$ws.Range("I63").Value = 4959.8
$ws.Range("K63").Value = 4959.8
$ws.Range("M63").Value = -4273.8
$ws.Range("H66").Value = 4959.8
$ws.Range("I66").Value = 4959.8
$ws.Range("K66").Value = 24799
$ws.Range("M66").Value = -21367
$ws.Range("H74").Value = 2504862.8
$ws.Range("I74").Value = 4170271.8
$ws.Range("K74").Value = 4170271.8
$ws.Range("M74").Value = -4169397.8
$ws.Range("H77").Value = 2504862.8
$ws.Range("I77").Value = 4170271.8
$ws.Range("K77").Value = 20851359
$ws.Range("M77").Value = -20846991
$ws.Range("H102").Value = 1984
$ws.Range("I102").Value = 1779.5333
$ws.Range("K102").Value = 1779.5333
$ws.Range("M102").Value = -157.5333000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1318.875
$ws.Range("I22").Value = 309.8
$ws.Range("K22").Value = 309.8
$ws.Range("M22").Value = -136.8
$ws.Range("H86").Value = 21075.5
$ws.Range("I86").Value = 36806.95
$ws.Range("J86").Value = 3493.2942
$ws.Range("K86").Value = 36806.95
$ws.Range("L86").Value = 3493.2942
$ws.Range("M86").Value = -35683.95
$ws.Range("N86").Value = -5739.2942
$ws.Range("H89").Value = 21075.5
$ws.Range("I89").Value = 36806.95
$ws.Range("J89").Value = 3493.2942
$ws.Range("K89").Value = 184034.75
$ws.Range("L89").Value = 17466.471
$ws.Range("M89").Value = -178418.75
$ws.Range("N89").Value = -28698.471
$ws.Range("H94").Value = 1487.7632
$ws.Range("I94").Value = 1079.8148
$ws.Range("J94").Value = 2489.0908
$ws.Range("K94").Value = 1079.8148
$ws.Range("L94").Value = 2489.0908
$ws.Range("M94").Value = -628.8148000000001
$ws.Range("N94").Value = -3391.0908

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1523.75
$ws.Range("J107").Value = 2330.7
$ws.Range("L107").Value = 2330.7
$ws.Range("N107").Value = -6170.7
$ws.Range("H132").Value = 3074.0527
$ws.Range("I132").Value = 2999.5881
$ws.Range("K132").Value = 8998.764299999999
$ws.Range("M132").Value = -6468.764299999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1273.25
$ws.Range("J86").Value = 1273.25
$ws.Range("L86").Value = 3819.75
$ws.Range("N86").Value = -6191.75
$ws.Range("H89").Value = 1273.25
$ws.Range("J89").Value = 1273.25
$ws.Range("L89").Value = 11459.25
$ws.Range("N89").Value = -23315.25
$ws.Range("H113").Value = 1244.4286
$ws.Range("I113").Value = 798.8333
$ws.Range("J113").Value = 1578.625
$ws.Range("K113").Value = 2396.4999
$ws.Range("L113").Value = 4735.875
$ws.Range("M113").Value = -226.4998999999998
$ws.Range("N113").Value = -9075.875
$ws.Range("H131").Value = 4101.1113
$ws.Range("I131").Value = 2809.818
$ws.Range("J131").Value = 4669.28
$ws.Range("K131").Value = 8429.454000000002
$ws.Range("L131").Value = 14007.84
$ws.Range("M131").Value = -3389.454000000002
$ws.Range("N131").Value = -24087.84
$ws.Range("H134").Value = 5372507.5
$ws.Range("I134").Value = 9678294
$ws.Range("K134").Value = 29034882
$ws.Range("M134").Value = -29029812

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5674.024
$ws.Range("I70").Value = 6224.636
$ws.Range("J70").Value = 5478.645
$ws.Range("K70").Value = 6224.636
$ws.Range("L70").Value = 5478.645
$ws.Range("M70").Value = -5954.636
$ws.Range("N70").Value = -6018.645
$ws.Range("H73").Value = 5674.024
$ws.Range("I73").Value = 6224.636
$ws.Range("J73").Value = 5478.645
$ws.Range("K73").Value = 6224.636
$ws.Range("L73").Value = 5478.645
$ws.Range("M73").Value = -5288.636
$ws.Range("N73").Value = -7350.645
$ws.Range("H113").Value = 843862.1
$ws.Range("I113").Value = 2177.1333
$ws.Range("K113").Value = 2177.1333
$ws.Range("M113").Value = -7.133299999999963
$ws.Range("H126").Value = 2292.5
$ws.Range("I126").Value = 2278.2
$ws.Range("J126").Value = 2328.25
$ws.Range("K126").Value = 6834.599999999999
$ws.Range("L126").Value = 6984.75
$ws.Range("M126").Value = -4364.599999999999
$ws.Range("N126").Value = -11924.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11370.154
$ws.Range("J22").Value = 2744
$ws.Range("L22").Value = 2744
$ws.Range("N22").Value = -3334
$ws.Range("H27").Value = 11370.154
$ws.Range("J27").Value = 2744
$ws.Range("L27").Value = 2744
$ws.Range("N27").Value = -2958
$ws.Range("H61").Value = 5117
$ws.Range("I61").Value = 1693.6364
$ws.Range("K61").Value = 1693.6364
$ws.Range("M61").Value = -1491.6364
$ws.Range("H93").Value = 2529771.2
$ws.Range("I93").Value = 2262.2666
$ws.Range("J93").Value = 7945862
$ws.Range("K93").Value = 2262.2666
$ws.Range("L93").Value = 7945862
$ws.Range("M93").Value = -1014.2666
$ws.Range("N93").Value = -7948358
$ws.Range("H106").Value = 23370
$ws.Range("J106").Value = 23370
$ws.Range("L106").Value = 23370
$ws.Range("N106").Value = -25894
$ws.Range("H113").Value = 5117
$ws.Range("I113").Value = 1693.6364
$ws.Range("K113").Value = 1693.6364
$ws.Range("M113").Value = 476.3635999999999
$ws.Range("H132").Value = 4715.1665
$ws.Range("I132").Value = 2971.8333
$ws.Range("K132").Value = 8915.499899999999
$ws.Range("M132").Value = -6385.499899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2531.4375
$ws.Range("I107").Value = 1550.0834
$ws.Range("J107").Value = 3120.25
$ws.Range("K107").Value = 4650.2502
$ws.Range("L107").Value = 9360.75
$ws.Range("M107").Value = -2730.2502
$ws.Range("N107").Value = -13200.75
$ws.Range("H132").Value = 626077.9
$ws.Range("I132").Value = 938.5
$ws.Range("J132").Value = 2501496
$ws.Range("K132").Value = 2815.5
$ws.Range("L132").Value = 7504488
$ws.Range("M132").Value = -285.5
$ws.Range("N132").Value = -7509548
